# Update the referenced paper number from the published revision
# "P0130R2" to the in-progress draft "D0130R1".
#
# Word keeps the text as a single run today; after a real edit in the
# middle of it, Word naturally splits it into two runs separated by the
# "last edit location" (_GoBack) bookmark, which also moves here from
# wherever it used to be (names are unique, so re-adding it elsewhere
# removes the old one). We reproduce both effects explicitly.

$d = $word.ActiveDocument

# Locate "P0130R2" without letting Find do the replacement itself (that
# would just rewrite the whole run as one piece and lose the split).
# Keep the Range in a variable so the Find narrows *that* instance down
# to the match (re-reading $d.Content afterwards would just hand back
# the whole document again).
$rng = $d.Content
$found = $rng.Find.Execute("P0130R2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start
    $end = $rng.End

    # Replace the whole token's text first.
    $whole = $d.Range($start, $end)
    $whole.Text = "D0130R1"

    # Split the run at the point between "D0130" and "R1" (5 characters
    # in) by dropping the _GoBack bookmark there; Word always keeps this
    # bookmark unique, so adding it here removes it from its prior
    # location automatically.
    $splitPoint = $start + 5
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
